$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Lazer" column header in E1
$ws.Range("E1").Value = "Lazer"

# Fill new column E (rows 2-8) with "-" placeholder, matching other columns
$ws.Range("E2").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("E5").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("E8").Value = "-"

# Copy style from column D onto the new column E so it matches the alternating row colors
$ws.Range("D1:D8").Copy()
$ws.Range("E1:E8").PasteSpecial(-4122)

# Update Lucas Henrique's row (row 5): Esportes (B5) time changes, and Relaxamento (D5) gets a value
$ws.Range("B5").Value = "2 Horas e 0 Minutos"
$ws.Range("D5").Value = "1 Horas e 0 Minutos"
